$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value2 = 2
$ws.Range("B7").Value2 = 43523
$ws.Range("C7").Value2 = 0.51388888888888895
$ws.Range("D7").Value2 = 0.56944444444444442
$ws.Range("G7").Value2 = "Proge."
$ws.Range("H7").Value2 = "MVC EF"
$ws.Range("F7").Formula = "=(D7-C7)*24*60 - E7"

# Row 8
$ws.Range("A8").Value2 = 3
$ws.Range("B8").Value2 = 43523
$ws.Range("C8").Value2 = 0.83333333333333337
$ws.Range("D8").Value2 = 0.91666666666666663
$ws.Range("G8").Value2 = "Proge."
$ws.Range("H8").Value2 = "MVC EF"
$ws.Range("F8").Formula = "=(D8-C8)*24*60 - E8"

# Restore the General number format (and border) on F7/F8 that Excel's
# auto-format heuristic overwrote when the formula referenced time cells.
$ws.Range("F9").Copy()
$ws.Range("F7:F8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection, matching the recorded cursor position.
$ws.Range("H11").Select()
